$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'43.975.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +1.92%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.250.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +1.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'317.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -0.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'100.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +2.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -1.10%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  +0.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  -3.32%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'37.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +0.71%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.0833"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +0.90%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'7.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  -1.04%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  -1.52%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'2.594.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +1.25%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'14.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  +1.02%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'0.854"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  -0.78%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'2.253.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +1.31%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'43.892.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +1.92%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'13.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -3.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'0.0₃0981"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +1.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'6.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -1.45%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'65.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  +0.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Formula = "'  -4.19%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'233.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  -0.95%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  -5.92%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Formula = "'  +0.29%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +6.65%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'38.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  +5.36%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Formula = "'  -0.78%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'6.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  -4.02%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'161.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +3.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'20.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  -0.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'0.0842"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -2.43%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Formula = "'  +1.20%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Formula = "'LidoDAOToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Formula = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Formula = "'3.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  -6.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Formula = "'Kaspa"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Formula = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Formula = "'0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +7.19%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Formula = "'ARBITRUM"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Formula = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Formula = "'1.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +6.09%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Formula = "'  -1.86%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'16.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  +17.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  -0.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'4.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  -5.88%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.0315"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -1.10%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Formula = "'  +0.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'1.771.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +2.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.196"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -3.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'74.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  +0.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'5.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -1.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'80.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -3.82%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'103.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  +0.85%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'57.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'1.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +4.60%  "
$ws.Range("E51").Style = "Normal"
